$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row that only held the teacher name under "Docentes responsáveis:"
# (original row 13: B13/C13 = "6270264 - Juan Fernando Zapata Zapata").
# Deleting it shifts every row below it up by one.
$ws.Rows(13).Delete()

# After the shift, several cells now need their text corrected so the grid
# matches the published version exactly.

# Row 10 (Objetivos:) now shows the teacher's name instead of the long text.
$ws.Range("B10").Value = "6270264 - Juan Fernando Zapata Zapata"
$ws.Range("C10").Value = "6270264 - Juan Fernando Zapata Zapata"

# Row 13 (Programa resumido:) now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) now shows the activation date. Copy the existing
# text-formatted date cell (row 8) instead of assigning the string
# directly, so Excel does not auto-convert it into a date serial number.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# Row 18 (Método:) now shows the teacher's name again.
$ws.Range("B18").Value = "6270264 - Juan Fernando Zapata Zapata"
$ws.Range("C18").Value = "6270264 - Juan Fernando Zapata Zapata"

# Row 19 (Critério:) now shows the evaluation method text.
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Row 20 (Norma de recuperação:) now shows the passing-grade criterion.
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# Row 21 (Bibliografia:) now shows the recovery-grade rule.
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
